# Rewords the "Part N - ..." section-sheet titles to "Model N - ..."
# on the eight C++ Models section/divider slides, without disturbing
# the existing run formatting (sz=3600), the soft line break before
# it, or the surrounding "C++ Models" text.

$p = $ppt.ActivePresentation

$targetSlides = @(6, 14, 17, 21, 24, 27, 31, 35)

foreach ($slideIndex in $targetSlides) {
    $slide = $p.Slides.Item($slideIndex)
    $titleShape = $slide.Shapes.Item(1)
    $tr = $titleShape.TextFrame.TextRange

    $fullText = $tr.Text
    $idx = $fullText.IndexOf("Part ")

    if ($idx -ge 0) {
        # Select from the start of "Part ..." through the end of the
        # text range - this spans exactly the existing run, so setting
        # its .Text keeps the result as a single run (matching how the
        # rest of the run's formatting, e.g. sz="3600", is preserved).
        $runLength = $tr.Length - $idx
        $sub = $tr.Characters($idx + 1, $runLength)
        $sub.Text = $sub.Text -replace "^Part ", "Model "
    }
}
